# The commit adds one new weekly price record (row) to the "Betarraga"
# sheet. In the original workbook this corresponds to inserting a brand
# new row at sheet row 174 and shifting every following record down by
# one row (so the former row 174 becomes row 175, ..., former row 249
# becomes row 250).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 174; existing rows 174.. shift down.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new weekly record.
$ws.Cells.Item(174, 1).Value  = 10
$ws.Cells.Item(174, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(174, 3).Value  = "La Araucanía"
$ws.Cells.Item(174, 4).Value  = 44466
$ws.Cells.Item(174, 5).Value  = 9
$ws.Cells.Item(174, 6).Value  = 100114014
$ws.Cells.Item(174, 7).Value  = "Betarraga"
$ws.Cells.Item(174, 8).Value  = "Sin especificar"
$ws.Cells.Item(174, 9).Value  = "Primera"
$ws.Cells.Item(174, 10).Value = 20
$ws.Cells.Item(174, 11).Value = 8000
$ws.Cells.Item(174, 12).Value = 8000
$ws.Cells.Item(174, 13).Value = 8000
$ws.Cells.Item(174, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(174, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(174, 16).Value = 667
$ws.Cells.Item(174, 17).Value = 12
$ws.Cells.Item(174, 18).Value = "Hortaliza"
